# Commit: "Doing Updates for Financials"
# Insert two new quarter columns (new D, E) before the existing data block
# (old D:K), shifting the existing 8 quarters of data right to F:M, and
# populate the two new columns with the newest quarter figures. A handful
# of cells in the previously-existing columns also carry corrected values
# (visible as real edits once re-aligned to their new column letters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns at D:E -- this shifts the old D:K data block
#    to F:M and leaves D:E blank with default formatting.
$ws.Columns("D:E").Insert()

# 2) Copy number formats (date format in header rows, number format in
#    data rows) from column F (the old column D, now shifted right) into
#    the new D:E columns so the new cells render the same way as their
#    neighbours. Restricted to the three contiguous data blocks (rows
#    5/6/37/79 are section-label-only rows with no D:K data and must stay
#    untouched).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rowData = @{
    "7" = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    "8" = @(1801000, 1550000, 1537000, 1517000, 1782000, 1663000, 1593000, 1478000, 1802000, 1677000)
    "9" = @(1359000, 1140000, 1134000, 1097000, 1276000, 1190000, 1133000, 1065000, 1323000, 1200000)
    "10" = @(442000, 410000, 403000, 420000, 506000, 473000, 460000, 413000, 479000, 477000)
    "11" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "12" = @(60000, 56000, 62000, 64000, 59000, 52000, 57000, 67000, 83000, 56000)
    "13" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "14" = @(80000, 9000, 209000, 12000, 6000, 4000, 4000, 5000, 3000, 10000)
    "15" = @(16000, 14000, 16000, 16000, 16000, 17000, 16000, 16000, 15000, 17000)
    "16" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "17" = @(1738000, 1425000, 1643000, 1408000, 1586000, 1463000, 1420000, 1361000, 1656000, 1488000)
    "18" = @(63000, 125000, -106000, 109000, 196000, 200000, 173000, 117000, 146000, 189000)
    "19" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "20" = @(40000, -10000, -9000, -5000, -24000, -8000, -2000, -7000, -17000, -8000)
    "21" = @(184000, 193000, -30000, 190000, 263000, 283000, 258000, 195000, 214000, 265000)
    "22" = @(43000, 43000, 41000, 41000, 41000, 42000, 41000, 39000, 40000, 41000)
    "23" = @(60000, 72000, -156000, 63000, 131000, 150000, 130000, 71000, 89000, 140000)
    "24" = @(56000, -15000, -12000, 7000, 34000, 31000, 33000, 14000, 17000, 31000)
    "25" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "26" = @(4000, 87000, -144000, 56000, 97000, 119000, 97000, 57000, 72000, 109000)
    "27" = @(-9000, 73000, -155000, 43000, 84000, 106000, 85000, -17000, 56000, 94000)
    "28" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "29" = @(-51000, -1000, -2000, -35000, -140000, 0, 5000, 0, -11000, -2000)
    "30" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "31" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "32" = @(-40000, 10000, 9000, 5000, 24000, 8000, 2000, 7000, 17000, 8000)
    "33" = @(-60000, 72000, -157000, 8000, -56000, 106000, 90000, -17000, 45000, 92000)
    "34" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "35" = @(-60000, 72000, -157000, 8000, -56000, 106000, 90000, -17000, 45000, 92000)
    "38" = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    "39" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "40" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "41" = @(464000, 334000, 343000, 348000, 537000, 405000, 377000, 401000, 498000, 318000)
    "42" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "43" = @(1356000, 1309000, 1272000, 1338000, 1270000, 1408000, 1321000, 1298000, 1282000, 1387000)
    "44" = @(806000, 916000, 842000, 822000, 780000, 824000, 828000, 800000, 699000, 776000)
    "45" = @(397000, 296000, 282000, 283000, 243000, 263000, 290000, 281000, 278000, 270000)
    "46" = @(3023000, 2855000, 2739000, 2791000, 2830000, 2900000, 2816000, 2780000, 2757000, 2751000)
    "47" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "48" = @(359000, 348000, 326000, 338000, 341000, 321000, 304000, 288000, 287000, 289000)
    "49" = @(3287000, 3090000, 3107000, 3302000, 3319000, 3332000, 3354000, 3377000, 3399000, 3441000)
    "50" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "51" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "52" = @(1092000, 1212000, 1192000, 1210000, 1164000, 1297000, 1293000, 1280000, 1230000, 1230000)
    "53" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "54" = @(7761000, 7505000, 7364000, 7641000, 7654000, 7850000, 7767000, 7725000, 7673000, 7711000)
    "55" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "56" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "57" = @(897000, 804000, 711000, 697000, 762000, 720000, 731000, 765000, 781000, 718000)
    "58" = @(185000, 246000, 160000, 34000, 52000, 269000, 267000, 252000, 50000, 256000)
    "59" = @(1200000, 954000, 1008000, 1113000, 1075000, 1057000, 1115000, 1142000, 1134000, 1028000)
    "60" = @(2282000, 2004000, 1879000, 1844000, 1889000, 2046000, 2113000, 2159000, 1965000, 2002000)
    "61" = @(2980000, 2881000, 2952000, 3038000, 2939000, 2984000, 3015000, 3076000, 3001000, 3033000)
    "62" = @(1227000, 1312000, 1313000, 1319000, 1279000, 1233000, 1192000, 1165000, 1146000, 1156000)
    "63" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "64" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "65" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "66" = @(6507000, 6215000, 6159000, 6220000, 6125000, 6281000, 6339000, 6419000, 6131000, 6204000)
    "67" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "68" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "69" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "70" = @(859000, 846000, 834000, 822000, 810000, 799000, 786000, 776000, 847000, 835000)
    "71" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "72" = @(606000, 666000, 594000, 782000, 857000, 913000, 806000, 716000, 867000, 822000)
    "73" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "74" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "75" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "76" = @(395000, 444000, 371000, 599000, 719000, 770000, 642000, 530000, 695000, 672000)
    "77" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "80" = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    "81" = @(-60000, 72000, -157000, 8000, -56000, 106000, 90000, -17000, 45000, 92000)
    "82" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "83" = @(81000, 78000, 85000, 86000, 91000, 91000, 87000, 85000, 85000, 84000)
    "84" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "85" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "86" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "87" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "88" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "89" = @(396000, 56000, 112000, -28000, 490000, 124000, 93000, 40000, 516000, 215000)
    "90" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "91" = @(-39000, -34000, -41000, -29000, -47000, -38000, -32000, -104000, -67000, -62000)
    "92" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "93" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "94" = @(-285000, -76000, -85000, -74000, -90000, -73000, -74000, -53000, -68000, -62000)
    "95" = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    "96" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "97" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "98" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "99" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "100" = @(35000, 17000, -19000, -91000, -268000, -27000, -47000, -92000, -244000, -167000)
    "101" = @(-13000, -4000, -13000, 5000, 0, 4000, 4000, 8000, -24000, 0)
    "102" = @(133000, -7000, -5000, -188000, 132000, 28000, -24000, -97000, 180000, -14000)
}


# 3) Write the final values for every data row, columns D through M
#    (D,E = the two new quarters; F..M = the existing 8 quarters, with a
#    few cells carrying corrected figures relative to a pure shift).
foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 4 + $i   # column D = 4
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item([int]$r, $col).Value = $v
        }
    }
}
